$wb = $excel.ActiveWorkbook

# Layer0 sheet (sheet1.xml)
$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.4514317921401603
$ws1.Range("C2").Value = -0.7792597947380407
$ws1.Range("B3").Value = 0.9497641650128565
$ws1.Range("C3").Value = -1.355450823862738
$ws1.Range("B4").Value = 0.1605528397728489
$ws1.Range("C4").Value = 0.3515043829056255

# Layer1 sheet (sheet2.xml)
$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -1.096488982880605
$ws2.Range("C2").Value = -0.3691835905611328
$ws2.Range("B3").Value = 0.7186302102988767
$ws2.Range("C3").Value = 0.247179417018035
$ws2.Range("B4").Value = -1.209618343029855
$ws2.Range("C4").Value = 1.01223598979142
